{"js": "// Apply the built-in \"List Bullet\" paragraph style to the document's\n// (only) paragraph. In the canonical OOXML this is the addition of\n//   <w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>\n// as the first child of the body's <w:p> \u2014 Word materializes the\n// \"List Bullet\" style (plus its backing numbering definition) into\n// styles.xml / numbering.xml the first time the style is used.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  const firstParagraph = paragraphs.items[0];\n  firstParagraph.style = \"List Bullet\";\n  await context.sync();\n}\n", "ps1": "# Apply the built-in \"List Bullet\" paragraph style to the document's\n# (only) paragraph. In the canonical OOXML this is the addition of\n#   <w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>\n# as the first child of the body's <w:p> -- Word materializes the\n# \"List Bullet\" style (plus its backing numbering definition) into\n# styles.xml / numbering.xml the first time the style is used.\n\n$d = $word.ActiveDocument\n\nif ($d.Paragraphs.Count -ge 1) {\n    $p = $d.Paragraphs(1)\n    $p.Range.Style = \"List Bullet\"\n}\n"}
